# Insert a new "Industry" column at column C, shifting Mutual Fund,
# Status, Jan_2026, Dec_2025, Oct_2025, MoM, QoQ columns one position to
# the right (this is the output of the motilal_portfolio_change_engine
# run which now also reports each holding's industry classification).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts C:I -> D:J)
$ws.Columns.Item(3).Insert()

# Header for the newly inserted column, matching style of neighboring headers
$ws.Cells.Item(1, 4).Copy()
$ws.Cells.Item(1, 3).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(1, 3).Value = "Industry"

# Populate the new column's data rows with "N.A." (industry unknown for these CDs/CPs)
$ws.Cells.Item(2, 3).Value = "N.A."
$ws.Cells.Item(3, 3).Value = "N.A."
$ws.Cells.Item(4, 3).Value = "N.A."
